$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(-0.08343332852268287, 0.2276087606115918, 0.08273643786598961, 0.2876394233515107, 0.06006953017233075, 22)
    3 = @(-0.07650533591609261, 0.2496604832070936, 0.09152595977496505, 0.3025325763863539, 0.06387222861430945, 22)
    4 = @(-0.126395414308652, 0.3014333265446221, 0.1564314459343321, 0.3955141538988612, 0.08178241899090313, 22)
    5 = @(-0.07996933221938775, 0.2147519892694937, 0.0649744043265079, 0.254900773491388, 0.05281563180892435, 22)
    6 = @(-0.1049143714156674, 0.2299523683858045, 0.0944954094530339, 0.3074010563629115, 0.06305265529074708, 22)
    7 = @(-0.1014503751123723, 0.2340950556680925, 0.09073056349745118, 0.3012151448673376, 0.06189021958543866, 22)
    8 = @(-0.09244185181752758, 0.212635985751665, 0.07142287205047154, 0.2672505791396373, 0.05471892511606922, 22)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
